# Update cryptocurrency price/volume data cells to match latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.652.12"
$ws.Range("E2").Value = "  -2.47%  "
$ws.Range("D3").Value = "2.006.32"
$ws.Range("E3").Value = "  -4.91%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  +0.62%  "
$ws.Range("D5").Value = "'332.00"
$ws.Range("E5").Value = "  -3.76%  "
$ws.Range("D6").Value = "'1.012"
$ws.Range("E6").Value = "  +0.59%  "
$ws.Range("D7").Value = "'0.5033"
$ws.Range("E7").Value = "  -3.81%  "
$ws.Range("D8").Value = "'0.4267"
$ws.Range("E8").Value = "  -4.09%  "
$ws.Range("D9").Value = "'54.88"
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D10").Value = "'0.09191"
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("D11").Value = "'1.128"
$ws.Range("E11").Value = "  -3.91%  "
$ws.Range("D12").Value = "'23.56"
$ws.Range("E12").Value = "  -5.74%  "
$ws.Range("D13").Value = "'8.148"
$ws.Range("E13").Value = "  -6.15%  "
$ws.Range("D14").Value = "2.005.79"
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").Value = "'6.556"
$ws.Range("E15").Value = "  -5.88%  "
$ws.Range("D16").Value = "'95.39"
$ws.Range("E16").Value = "  -6.44%  "
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "'0.00001125"
$ws.Range("E18").Value = "  -3.36%  "
$ws.Range("D19").Value = "'0.06683"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("D20").Value = "'19.93"
$ws.Range("E20").Value = "  -6.09%  "
$ws.Range("D21").Value = "'1.011"
$ws.Range("E21").Value = "  +0.54%  "
$ws.Range("D22").Value = "'5.991"
$ws.Range("E22").Value = "  -5.64%  "
$ws.Range("D23").Value = "29.646.27"
$ws.Range("E23").Value = "  -2.57%  "
$ws.Range("D24").Value = "'12.08"
$ws.Range("E24").Value = "  -4.52%  "
$ws.Range("D25").Value = "'2.282"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").Value = "'159.23"
$ws.Range("E26").Value = "  -2.43%  "
$ws.Range("D27").Value = "'20.84"
$ws.Range("E27").Value = "  -5.61%  "
$ws.Range("D28").Value = "'6.433"
$ws.Range("E28").Value = "  -5.74%  "
$ws.Range("D29").Value = "'2.334"
$ws.Range("E29").Value = "  -8.15%  "
$ws.Range("D30").Value = "'129.07"
$ws.Range("E30").Value = "  -3.81%  "
$ws.Range("D31").Value = "'1.067"
$ws.Range("E31").Value = "  -7.65%  "
$ws.Range("D32").Value = "'1.584"
$ws.Range("E32").Value = "  -8.88%  "
$ws.Range("D33").Value = "'0.09967"
$ws.Range("E33").Value = "  -5.57%  "
$ws.Range("D34").Value = "'5.867"
$ws.Range("E34").Value = "  -6.53%  "
$ws.Range("D35").Value = "'3.807"
$ws.Range("E35").Value = "  -2.84%  "
$ws.Range("D36").Value = "'9.581"
$ws.Range("E36").Value = "  -7.92%  "
$ws.Range("D37").Value = "'0.02482"
$ws.Range("E37").Value = "  -5.42%  "
$ws.Range("D38").Value = "'1.322"
$ws.Range("E38").Value = "  -2.74%  "
$ws.Range("D39").Value = "'0.06398"
$ws.Range("E39").Value = "  -5.79%  "
$ws.Range("D40").Value = "'0.6607"
$ws.Range("E40").Value = "  -6.43%  "
$ws.Range("D41").Value = "'11.79"
$ws.Range("E41").Value = "  -6.34%  "
$ws.Range("D42").Value = "'0.2076"
$ws.Range("E42").Value = "  -6.79%  "
$ws.Range("D43").Value = "'1.011"
$ws.Range("E43").Value = "  +0.48%  "
$ws.Range("D44").Value = "'0.6379"
$ws.Range("E44").Value = "  -7.06%  "
$ws.Range("D45").Value = "'13.67"
$ws.Range("E45").Value = "  -5.87%  "
$ws.Range("D46").Value = "'2.220"
$ws.Range("E46").Value = "  -6.03%  "
$ws.Range("E47").Value = "  -4.83%  "
$ws.Range("D48").Value = "'3.534"
$ws.Range("E48").Value = "  -3.02%  "
$ws.Range("D49").Value = "'0.07009"
$ws.Range("D50").Value = "'0.00000000326"
$ws.Range("E50").Value = "  -5.59%  "
$ws.Range("D51").Value = "'1.134"
$ws.Range("E51").Value = "  -6.05%  "
